$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "dlgProfil"
$ws.Range("B5").Value = "dlgBaseNormal"
$ws.Range("D5").Value = "Mein Profil | TT-Planer"

$ws.Range("E13").Select()
